# Update countries & provincias Spain
# - Update "Estados Unidos" row (row 4) totals
# - Re-rank a few country rows whose values changed so their order (and thus
#   the shared-string / row ordering) swaps between each pair:
#     Santa Lucia (row 188)   <-> Belice (row 189)
#     Namibia     (row 194)   <-> San Vicente y las Granadinas (row 195)
#     Burundi     (row 198)   <-> San Cristobal y Nieves (row 199)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update Estados Unidos row (row 4) figures
$ws.Range("B4").Value = 1188385
$ws.Range("C4").Value = 263
$ws.Range("E4").Value = 941191

# 2) Swap the full data rows (columns A:H) for the three country pairs whose
#    ranking order changed.
function Swap-Rows($sheet, $r1, $r2) {
    $rng1 = $sheet.Range("A$r1`:H$r1")
    $rng2 = $sheet.Range("A$r2`:H$r2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

Swap-Rows $ws 188 189
Swap-Rows $ws 194 195
Swap-Rows $ws 198 199
